$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "total price per item" column
$ws.Range("Q1").Value = "Eszköz össz ár"

# Per-row totals: PRODUCT(price)*(qty), formatted as currency (Ft)
$ws.Range("Q2:Q12").NumberFormat = '#,##0\ "Ft"'
for ($r = 2; $r -le 12; $r++) {
    $ws.Range("Q$r").Formula = "=PRODUCT(C$r)*(D$r)"
}

# Keep the same number format going down to the blank row under the table
$ws.Range("Q14").NumberFormat = '#,##0\ "Ft"'

# Grand total row
$ws.Range("O15").Value = "Összköltség:"
$ws.Range("Q15").NumberFormat = '#,##0\ "Ft"'
$ws.Range("Q15").Formula = "=SUM(Q2:Q12)"

# Widen the new column and select the grand-total label cell, like in the source file
$ws.Columns.Item(17).ColumnWidth = 17.75
[void]$ws.Range("O15").Select()
